# Add season-record columns (Wins / Losses / Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1) onto
# the three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-51) gets the team's season record repeated across
# the new columns.
for ($r = 2; $r -le 51; $r++) {
    $ws.Range("AD$r").Value = 64
    $ws.Range("AE$r").Value = 98
    $ws.Range("AF$r").Value = 0
}
